# Add two new columns: I ("I0") and J ("IF").
# Header style mirrors the existing header cells (bold, bordered, centered),
# so copy H1's format into I1:J1 before setting the new header text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I column is always 1; J column mirrors the existing IP (H) column value.
for ($r = 2; $r -le 39; $r++) {
    $ipValue = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $ipValue
}
